$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F-column (time_taken) timestamps on the "data" sheet ---
# These reflect a later re-run of the panel-scraping script (rows 2..99).
$newTimes = @(
    "2021-10-05 14:22:24.408650",
    "2021-10-05 14:22:24.408657",
    "2021-10-05 14:22:24.408659",
    "2021-10-05 14:22:24.408661",
    "2021-10-05 14:22:24.408663",
    "2021-10-05 14:22:24.408665",
    "2021-10-05 14:22:24.408667",
    "2021-10-05 14:22:24.408669",
    "2021-10-05 14:22:24.408671",
    "2021-10-05 14:22:24.408673",
    "2021-10-05 14:22:24.408675",
    "2021-10-05 14:22:24.408677",
    "2021-10-05 14:22:24.408679",
    "2021-10-05 14:22:24.408681",
    "2021-10-05 14:22:24.408683",
    "2021-10-05 14:22:24.408685",
    "2021-10-05 14:22:24.408687",
    "2021-10-05 14:22:24.408689",
    "2021-10-05 14:22:24.408691",
    "2021-10-05 14:22:24.408694",
    "2021-10-05 14:22:24.408696",
    "2021-10-05 14:22:24.408698",
    "2021-10-05 14:22:24.408699",
    "2021-10-05 14:22:24.408701",
    "2021-10-05 14:22:24.408704",
    "2021-10-05 14:22:24.408706",
    "2021-10-05 14:22:24.408708",
    "2021-10-05 14:22:24.408710",
    "2021-10-05 14:22:24.408712",
    "2021-10-05 14:22:24.408714",
    "2021-10-05 14:22:24.408716",
    "2021-10-05 14:22:24.408718",
    "2021-10-05 14:22:24.408720",
    "2021-10-05 14:22:24.408723",
    "2021-10-05 14:22:24.408725",
    "2021-10-05 14:22:24.408726",
    "2021-10-05 14:22:24.408728",
    "2021-10-05 14:22:24.408730",
    "2021-10-05 14:22:24.408732",
    "2021-10-05 14:22:24.408734",
    "2021-10-05 14:22:24.408736",
    "2021-10-05 14:22:24.408738",
    "2021-10-05 14:22:24.408740",
    "2021-10-05 14:22:24.408742",
    "2021-10-05 14:22:24.408744",
    "2021-10-05 14:22:24.408746",
    "2021-10-05 14:22:24.408748",
    "2021-10-05 14:22:24.408750",
    "2021-10-05 14:22:24.408752",
    "2021-10-05 14:22:24.408754",
    "2021-10-05 14:22:24.408756",
    "2021-10-05 14:22:24.408758",
    "2021-10-05 14:22:24.408760",
    "2021-10-05 14:22:24.408762",
    "2021-10-05 14:22:24.408764",
    "2021-10-05 14:22:24.408766",
    "2021-10-05 14:22:24.408768",
    "2021-10-05 14:22:24.408770",
    "2021-10-05 14:22:24.408772",
    "2021-10-05 14:22:24.408774",
    "2021-10-05 14:22:24.408776",
    "2021-10-05 14:22:24.408778",
    "2021-10-05 14:22:24.408781",
    "2021-10-05 14:22:24.408783",
    "2021-10-05 14:22:24.408786",
    "2021-10-05 14:22:24.408788",
    "2021-10-05 14:22:24.408791",
    "2021-10-05 14:22:24.408793",
    "2021-10-05 14:22:24.408795",
    "2021-10-05 14:22:24.408797",
    "2021-10-05 14:22:24.408799",
    "2021-10-05 14:22:24.408801",
    "2021-10-05 14:22:24.408803",
    "2021-10-05 14:22:24.408806",
    "2021-10-05 14:22:24.408808",
    "2021-10-05 14:22:24.408810",
    "2021-10-05 14:22:24.408814",
    "2021-10-05 14:22:24.408816",
    "2021-10-05 14:22:24.408818",
    "2021-10-05 14:22:24.408820",
    "2021-10-05 14:22:24.408823",
    "2021-10-05 14:22:24.408825",
    "2021-10-05 14:22:24.408827",
    "2021-10-05 14:22:24.408829",
    "2021-10-05 14:22:24.408831",
    "2021-10-05 14:22:24.408833",
    "2021-10-05 14:22:24.408836",
    "2021-10-05 14:22:24.408838",
    "2021-10-05 14:22:24.408840",
    "2021-10-05 14:22:24.408842",
    "2021-10-05 14:22:24.408844",
    "2021-10-05 14:22:24.408846",
    "2021-10-05 14:22:24.408850",
    "2021-10-05 14:22:24.408852",
    "2021-10-05 14:22:24.408854",
    "2021-10-05 14:22:24.408857",
    "2021-10-05 14:22:24.408859",
    "2021-10-05 14:22:24.408861"
)

$row = 2
foreach ($t in $newTimes) {
    $dataSheet.Range("F$row").Value = $t
    $row++
}

# --- Add a new "metadata" sheet after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Match formatting used on the "data" sheet: bold/centered/bordered header
# row (B1:G1) and the bold/centered/bordered index cell (A2).
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Rare anaemia"
$metaSheet.Range("C2").Value = 518

# data_version "1.27" must stay textual, not become the number 1.27.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.27"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-10-01T08:46:36.144101Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:24.406202"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/518/?format=json"

# Keep "data" as the active/selected sheet (unchanged by this edit).
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null

Write-Host "metadata sheet added; F column timestamps refreshed"
